# Regenerate save_data column G ("K") values to reflect the new K-based
# (strikeouts) computation instead of the old Strike# based one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newValues = @{
    2  = 6
    3  = 3
    4  = 7
    5  = 2
    6  = 1
    7  = 0
    9  = 5
    10 = 7
    11 = 6
    12 = 5
    13 = 2
    14 = 1
    15 = 6
    16 = 1
    17 = 2
    18 = 2
    19 = 3
    20 = 1
    21 = 5
    22 = 4
    23 = 5
    24 = 6
    25 = 4
    26 = 6
    27 = 2
    28 = 4
    29 = 8
    30 = 8
    31 = 2
    32 = 3
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
